$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set A9 to the next date in the sequence (continuing from A8 = 41207)
# and B9 to 1 hour, matching the style already applied to the row.
$ws.Range("A9").Value = 41208
$ws.Range("B9").Value = 1

# Recalculate formulas (D5 = SUM(B5:B24), F5 = D5*E5) so cached values update.
$excel.CalculateFullRebuild()

# Update the active selection to A10, as shown in the sheetView.
$ws.Range("A10").Select()
